# Rename the worksheet tab from "Neurology" to "Session"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Session"

# Remove the trailing 4 log rows (rows 37-40), which shifts the sheet's
# used range/dimension from A1:F40 down to A1:F36.
$ws.Rows("37:40").Delete()
